$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1,2,4,1,6,1,3,1,0,2,0,6,8,6,3,6,7,5,2,1,4)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
